# Progress update as of 04-Nov-2025:
# On the "Training Dashboard" sheet, for every data row (3-27):
#   - PERIOD TO EXPIRE (column H) decreases by 1 day
#   - LAST UPDATE (column I) moves from 03-Nov-2025 to 04-Nov-2025
#
# Note: plain `.Value` assignment of a date-shaped string such as
# "04-Nov-2025" makes Excel auto-convert the cell to a real date
# (changing its number format/style). To keep the cell as literal text
# with its original style intact, we build the text via a formula and
# then convert that formula to a static value with Paste Special
# (values only), which does not trigger date re-interpretation.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Training Dashboard")

for ($row = 3; $row -le 27; $row++) {
    $hCell = $ws.Cells.Item($row, 8)   # column H: PERIOD TO EXPIRE
    $hCell.Value2 = $hCell.Value2 - 1

    $iCell = $ws.Cells.Item($row, 9)   # column I: LAST UPDATE
    $iCell.Formula = "=""04-Nov-2025"""
    $iCell.Copy() | Out-Null
    $iCell.PasteSpecial(-4163)         # xlPasteValues
}

$excel.CutCopyMode = $false
